# Extend the example variable register (inst/extdata/example_register.xlsx)
# with a row for a pooled *ordinal* variable, so the example data set
# covers all variable types (metric / ordinal / categorical, pooled and
# not pooled).
#
# The existing row 9, which used to describe "pvkat_pooled" (a pooled
# categorical variable), now describes "pvord_pooled" (a pooled ordinal
# variable); a brand-new row 10 is appended describing "pvkat_pooled".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: relabel from "pvkat_pooled" to "pvord_pooled" (Nr/kw1/kw2 unchanged)
$ws.Range("B9").Value = "pvord_pooled"

# Row 10 (new): Nr = 9, varName = "pvkat_pooled", no keywords set
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "pvkat_pooled"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Style = "Normal"
